# The document contains a single MERGEFIELD-like field whose code reads:
#   m:'contents.txt'.fromMediaWikiURI()
# stored as a real Word field (fldChar begin/end + instrText runs, with a
# "_GoBack" bookmark nested between "MediaWiki" and "URI"). The commit
# replaces that field with plain literal text runs -- "{" + the same
# tokens + "}" -- using <w:t> instead of <w:instrText>, dropping the
# fldChar wrapper and the two pure-whitespace instrText runs that used to
# sit right after "begin" and right before "end". The bookmark stays put,
# between "MediaWiki" and "URI".

$d = $word.ActiveDocument

# Find the paragraph that owns the field (there is exactly one field in
# this document). Locate it by scanning $d.Paragraphs for the one whose
# range contains the field's code -- going through Field.Code.Paragraphs
# / Field.Result.Paragraphs directly is unreliable when the field result
# is empty.
$field = $d.Fields.Item(1)
$codeStart = $field.Code.Start
$fieldPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($codeStart -ge $candidate.Range.Start -and $codeStart -lt $candidate.Range.End) {
        $fieldPara = $candidate
    }
}

# Use the paragraph's own (live) Range property directly -- rebuilding an
# equivalent range via $d.Range(start, end) does not reliably span a
# field whose result is empty (it collapses to an insertion point instead
# of replacing the whole paragraph).
$target = $fieldPara.Range

# Build the replacement paragraph: one run per literal token (mirroring
# the original run splitting), with the "_GoBack" bookmark re-inserted
# between "MediaWiki" and "URI".
$tokens = @("{", "m", ":", "'", "contents.txt", "'", ".from", "MediaWiki")
$xml = ""
foreach ($tok in $tokens) {
    $xml += "<w:r><w:t>" + $tok + "</w:t></w:r>"
}
$xml += '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
$xml += "<w:r><w:t>URI</w:t></w:r>"
$xml += "<w:r><w:t>()</w:t></w:r>"
$xml += '<w:r><w:t xml:space="preserve">}</w:t></w:r>'

[void]$target.InsertXML("<w:p>" + $xml + "</w:p>")
